$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "ECs"
$ws.Range("H2").Value = 0.6219589999999999
$ws.Range("I2").Value = 0.6398583988494134
$ws.Range("J2").Value = 0.6398583988494134
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.5706193333333334
$ws.Range("N2").Value = 1.711858
$ws.Range("O2").Value = 0.4188640502130462
$ws.Range("P2").Value = 0.4188640502130463
$ws.Range("Q2").Value = 0.1183006099802222
$ws.Range("R2").Value = 1.064705489822
$ws.Range("S2").Value = 0.2680136805049
$ws.Range("T2").Value = 0.2680136805049001
# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("H3").Value = 0.6219589999999999
$ws.Range("I3").Value = 0.6398583988494134
$ws.Range("J3").Value = 0.6398583988494134
$ws.Range("M3").Value = 0.4846943333333333
$ws.Range("N3").Value = 1.454083
$ws.Range("O3").Value = 0.3557906641356566
$ws.Range("P3").Value = 0.3557906641356566
$ws.Range("Q3").Value = 0.1004866676218889
$ws.Range("R3").Value = 0.9043800085969997
$ws.Range("S3").Value = 0.2276556446794106
$ws.Range("T3").Value = 0.2276556446794107
# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("H4").Value = 0.6219589999999999
$ws.Range("I4").Value = 0.6398583988494134
$ws.Range("J4").Value = 0.6398583988494134
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3069883333333334
$ws.Range("N4").Value = 0.920965
$ws.Range("O4").Value = 0.2253452856512971
$ws.Range("P4").Value = 0.2253452856512971
$ws.Range("Q4").Value = 0.06364471893722222
$ws.Range("R4").Value = 0.572802470435
$ws.Range("S4").Value = 0.1441890736651027
$ws.Range("T4").Value = 0.1441890736651027
# Row 5
$ws.Range("D5").Value = "ECs"
$ws.Range("G5").Value = 0.116689
$ws.Range("H5").Value = 0.350067
$ws.Range("I5").Value = 0.3601416011505865
$ws.Range("J5").Value = 0.3601416011505865
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.5706193333333334
$ws.Range("N5").Value = 1.711858
$ws.Range("O5").Value = 0.4188640502130462
$ws.Range("P5").Value = 0.4188640502130463
$ws.Range("Q5").Value = 0.06658499938733334
$ws.Range("R5").Value = 0.599264994486
$ws.Range("S5").Value = 0.1508503697081461
$ws.Range("T5").Value = 0.1508503697081462
# Row 6
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 0.116689
$ws.Range("H6").Value = 0.350067
$ws.Range("I6").Value = 0.3601416011505865
$ws.Range("J6").Value = 0.3601416011505865
$ws.Range("M6").Value = 0.4846943333333333
$ws.Range("N6").Value = 1.454083
$ws.Range("O6").Value = 0.3557906641356566
$ws.Range("P6").Value = 0.3557906641356566
$ws.Range("Q6").Value = 0.05655849706233333
$ws.Range("R6").Value = 0.509026473561
$ws.Range("S6").Value = 0.1281350194562459
$ws.Range("T6").Value = 0.1281350194562459
# Row 7
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 0.116689
$ws.Range("H7").Value = 0.350067
$ws.Range("I7").Value = 0.3601416011505865
$ws.Range("J7").Value = 0.3601416011505865
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3069883333333334
$ws.Range("N7").Value = 0.920965
$ws.Range("O7").Value = 0.2253452856512971
$ws.Range("P7").Value = 0.2253452856512971
$ws.Range("Q7").Value = 0.03582216162833333
$ws.Range("R7").Value = 0.322399454655
$ws.Range("S7").Value = 0.08115621198619442
$ws.Range("T7").Value = 0.08115621198619442
